$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 16.92782909523858
$ws.Range("C2").Value = 7.871687853525354
$ws.Range("D2").Value = 14.11080947616109
$ws.Range("E2").Value = 14.70496088750325
$ws.Range("G2").Value = 3.742298591950208
$ws.Range("J2").Value = 8.790875988560462
$ws.Range("K2").Value = 12.59469174907704
$ws.Range("L2").Value = 12.22736327229394
$ws.Range("O2").Value = 36.0929145026763

$ws.Range("B3").Value = 16.79876423228349
$ws.Range("C3").Value = 7.854357191216713
$ws.Range("D3").Value = 14.10600833433263
$ws.Range("E3").Value = 14.72575277539706
$ws.Range("G3").Value = 3.744627015637752
$ws.Range("J3").Value = 8.802379848301793
$ws.Range("K3").Value = 12.50709371544188
$ws.Range("L3").Value = 12.23506672896719
$ws.Range("O3").Value = 36.1572123538184

$ws.Range("B4").Value = 16.72268098489688
$ws.Range("C4").Value = 7.843483876539338
$ws.Range("D4").Value = 14.10554292032484
$ws.Range("E4").Value = 14.7404468617251
$ws.Range("G4").Value = 3.746132803144115
$ws.Range("J4").Value = 8.809842448333043
$ws.Range("K4").Value = 12.45554804961841
$ws.Range("L4").Value = 12.24147353591535
$ws.Range("O4").Value = 36.20248435994559

$ws.Range("B5").Value = 16.69250194231486
$ws.Range("C5").Value = 7.83899504712196
$ws.Range("D5").Value = 14.10597899449291
$ws.Range("E5").Value = 14.74691990139023
$ws.Range("G5").Value = 3.746765628433608
$ws.Range("J5").Value = 8.812984196128172
$ws.Range("K5").Value = 12.43512554174196
$ws.Range("L5").Value = 12.2445067247442
$ws.Range("O5").Value = 36.22238822243487

$ws.Range("B6").Value = 16.68754144522325
$ws.Range("C6").Value = 7.838246200941315
$ws.Range("D6").Value = 14.106089234735
$ws.Range("E6").Value = 14.74802405082557
$ws.Range("G6").Value = 3.746871870333203
$ws.Range("J6").Value = 8.813511970462054
$ws.Range("K6").Value = 12.43177015007622
$ws.Range("L6").Value = 12.24503591195481
$ws.Range("O6").Value = 36.22578108346239

$ws.Range("B7").Value = 16.72227059864837
$ws.Range("C7").Value = 7.843423572210662
$ws.Range("D7").Value = 14.10554626604975
$ws.Range("E7").Value = 14.74053219493672
$ws.Range("G7").Value = 3.746141259806059
$ws.Range("J7").Value = 8.809884411003708
$ws.Range("K7").Value = 12.45527023982082
$ws.Range("L7").Value = 12.24151273156311
$ws.Range("O7").Value = 36.20274690071813

$ws.Range("B8").Value = 16.88268951772557
$ws.Range("C8").Value = 7.865760751085937
$ws.Range("D8").Value = 14.10863977518332
$ws.Range("E8").Value = 14.71173002650864
$ws.Range("G8").Value = 3.743085669975361
$ws.Range("J8").Value = 8.794759849396753
$ws.Range("K8").Value = 12.56403505674272
$ws.Range("L8").Value = 12.22967184016906
$ws.Range("O8").Value = 36.11388123502358

$ws.Range("B9").Value = 17.22098163926573
$ws.Range("C9").Value = 7.907699488013617
$ws.Range("D9").Value = 14.13432739308466
$ws.Range("E9").Value = 14.67052997717732
$ws.Range("G9").Value = 3.73769485763943
$ws.Range("J9").Value = 8.768254356749864
$ws.Range("K9").Value = 12.79419512527039
$ws.Range("L9").Value = 12.21972358867351
$ws.Range("O9").Value = 35.98564236215817

$ws.Range("B10").Value = 17.48204417441636
$ws.Range("C10").Value = 7.937347906227418
$ws.Range("D10").Value = 14.16504322700539
$ws.Range("E10").Value = 14.6495540449547
$ws.Range("G10").Value = 3.73409676277829
$ws.Range("J10").Value = 8.750684393751074
$ws.Range("K10").Value = 12.97232213130049
$ws.Range("L10").Value = 12.22045692317844
$ws.Range("O10").Value = 35.91956241484603

$ws.Range("B11").Value = 17.60309685027047
$ws.Range("C11").Value = 7.950577280849723
$ws.Range("D11").Value = 14.18155674008269
$ws.Range("E11").Value = 14.64202388037937
$ws.Range("G11").Value = 3.732537773426724
$ws.Range("J11").Value = 8.743100688437815
$ws.Range("K11").Value = 13.0550368923673
$ws.Range("L11").Value = 12.22252461749318
$ws.Range("O11").Value = 35.89562339008701

$ws.Range("B12").Value = 17.64922928697612
$ws.Range("C12").Value = 7.955549295319626
$ws.Range("D12").Value = 14.18817210390288
$ws.Range("E12").Value = 14.63946111216706
$ws.Range("G12").Value = 3.731958548476852
$ws.Range("J12").Value = 8.740287438008911
$ws.Range("K12").Value = 13.08657676856737
$ws.Range("L12").Value = 12.22355572811767
$ws.Range("O12").Value = 35.88743912331092

$ws.Range("B13").Value = 17.63928143538104
$ws.Range("C13").Value = 7.954480167403593
$ws.Range("D13").Value = 14.18673132031373
$ws.Range("E13").Value = 14.64000021799149
$ws.Range("G13").Value = 3.732082800811364
$ws.Range("J13").Value = 8.740890722833567
$ws.Range("K13").Value = 13.07977480797788
$ws.Range("L13").Value = 12.22332264569303
$ws.Range("O13").Value = 35.88916256345096

$ws.Range("B14").Value = 17.60688654870903
$ws.Range("C14").Value = 7.950987089635531
$ws.Range("D14").Value = 14.18209375393575
$ws.Range("E14").Value = 14.64180725655369
$ws.Range("G14").Value = 3.732489897497092
$ws.Range("J14").Value = 8.742868069078783
$ws.Range("K14").Value = 13.05762747761539
$ws.Range("L14").Value = 12.22260448238113
$ws.Range("O14").Value = 35.89493240727928

$ws.Range("B15").Value = 17.5870807097363
$ws.Range("C15").Value = 7.948842556509379
$ws.Range("D15").Value = 14.17930015759535
$ws.Range("E15").Value = 14.64295170464445
$ws.Range("G15").Value = 3.73274070379405
$ws.Range("J15").Value = 8.744086865545531
$ws.Range("K15").Value = 13.04408920213049
$ws.Range("L15").Value = 12.22219686048467
$ws.Range("O15").Value = 35.89858134253656

$ws.Range("B16").Value = 17.47417606369189
$ws.Range("C16").Value = 7.936478097001129
$ws.Range("D16").Value = 14.16401488121282
$ws.Range("E16").Value = 14.65008659642236
$ws.Range("G16").Value = 3.734200207061421
$ws.Range("J16").Value = 8.75118821258725
$ws.Range("K16").Value = 12.96694831894519
$ws.Range("L16").Value = 12.220356585241
$ws.Range("O16").Value = 35.9212501250645

$ws.Range("B17").Value = 17.40547372641625
$ws.Range("C17").Value = 7.928826445576584
$ws.Range("D17").Value = 14.15528635331348
$ws.Range("E17").Value = 14.65497853296937
$ws.Range("G17").Value = 3.735115451065
$ws.Range("J17").Value = 8.755649207781808
$ws.Range("K17").Value = 12.92003879457408
$ws.Range("L17").Value = 12.219670929929
$ws.Range("O17").Value = 35.93672503101035

$ws.Range("B18").Value = 17.36617615146292
$ws.Range("C18").Value = 7.924401183921427
$ws.Range("D18").Value = 14.15050529758358
$ws.Range("E18").Value = 14.65798165220072
$ws.Range("G18").Value = 3.735649201705691
$ws.Range("J18").Value = 8.758253563595531
$ws.Range("K18").Value = 12.89321766339185
$ws.Range("L18").Value = 12.21943986144407
$ws.Range("O18").Value = 35.94620190405283

$ws.Range("B19").Value = 17.35290927868404
$ws.Range("C19").Value = 7.922898721988787
$ws.Range("D19").Value = 14.14892772037307
$ws.Range("E19").Value = 14.65903100375058
$ws.Range("G19").Value = 3.735831180589828
$ws.Range("J19").Value = 8.759141976121427
$ws.Range("K19").Value = 12.88416471052176
$ws.Range("L19").Value = 12.21938970905714
$ws.Range("O19").Value = 35.94950953318509

$ws.Range("B20").Value = 17.4127648958073
$ws.Range("C20").Value = 7.929643489324373
$ws.Range("D20").Value = 14.1561907697737
$ws.Range("E20").Value = 14.65443817898643
$ws.Range("G20").Value = 3.735017263929819
$ws.Range("J20").Value = 8.755170343641714
$ws.Range("K20").Value = 12.92501601250198
$ws.Range("L20").Value = 12.21972702659035
$ws.Range("O20").Value = 35.93501806894847

$ws.Range("B21").Value = 17.6163940760223
$ws.Range("C21").Value = 7.95201411806005
$ws.Range("D21").Value = 14.18344612137618
$ws.Range("E21").Value = 14.64126865398896
$ws.Range("G21").Value = 3.732370021709057
$ws.Range("J21").Value = 8.742285688123426
$ws.Range("K21").Value = 13.06412697692246
$ws.Range("L21").Value = 12.22280870087252
$ws.Range("O21").Value = 35.89321375176905

$ws.Range("B22").Value = 17.75116441023096
$ws.Range("C22").Value = 7.966414820318904
$ws.Range("D22").Value = 14.20336765135882
$ws.Range("E22").Value = 14.63434430729717
$ws.Range("G22").Value = 3.730704746555981
$ws.Range("J22").Value = 8.73420587757454
$ws.Range("K22").Value = 13.15630078540035
$ws.Range("L22").Value = 12.22626833429929
$ws.Range("O22").Value = 35.87102706524466

$ws.Range("B23").Value = 17.67909303303729
$ws.Range("C23").Value = 7.958749198949555
$ws.Range("D23").Value = 14.19254339095576
$ws.Range("E23").Value = 14.63788619804098
$ws.Range("G23").Value = 3.731587620242775
$ws.Range("J23").Value = 8.73848710763351
$ws.Range("K23").Value = 13.10699906528131
$ws.Range("L23").Value = 12.22429002933656
$ws.Range("O23").Value = 35.88239850335564

$ws.Range("B24").Value = 17.40946793131197
$ws.Range("C24").Value = 7.929274185495697
$ws.Range("D24").Value = 14.15578114434123
$ws.Range("E24").Value = 14.65468187903477
$ws.Range("G24").Value = 3.73506163076536
$ws.Range("J24").Value = 8.755386714528314
$ws.Range("K24").Value = 12.92276535008634
$ws.Range("L24").Value = 12.21970115708441
$ws.Range("O24").Value = 35.93578797937828

$ws.Range("B25").Value = 17.12712226111092
$ws.Range("C25").Value = 7.896556374003863
$ws.Range("D25").Value = 14.12528861082651
$ws.Range("E25").Value = 14.68004178752694
$ws.Range("G25").Value = 3.739089265329104
$ws.Range("J25").Value = 8.775089136423958
$ws.Range("K25").Value = 12.73025160565281
$ws.Range("L25").Value = 12.22099863616965
$ws.Range("O25").Value = 36.01539801130588

